$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Is this appointment for you USER? "
$ws.Range("B2").Value = "NAME_CHECK_P"
$ws.Range("C2").Value = "P"

# Row 3
$ws.Range("A3").Value = "So is this for you USER?"
$ws.Range("B3").Value = "NAME_CHECK_P"
$ws.Range("C3").Value = "P"

# Row 4
$ws.Range("A4").Value = "So the appointment for USER? Isn't it ?"
$ws.Range("B4").Value = "NAME_CHECK_P"
$ws.Range("C4").Value = "P"

# Row 5
$ws.Range("A5").Value = "I'm making the reservation for USER. Okay ? :)"
$ws.Range("B5").Value = "NAME_CHECK_P"
$ws.Range("C5").Value = "P"

# Row 6
$ws.Range("A6").Value = "So could you please tell me the person name for the appointment?"
$ws.Range("B6").Value = "NAME_CHECK_NA"
$ws.Range("C6").Value = "NA"

# Row 7
$ws.Range("A7").Value = "First I need the person name for the appointment :)"
$ws.Range("B7").Value = "NAME_CHECK_NA"
$ws.Range("C7").Value = "NA"

# Row 8
$ws.Range("A8").Value = "Sorry I'm unable find an appointment person name. Could you please tell me the name ?"
$ws.Range("B8").Value = "NAME_CHECK_ERROR"
$ws.Range("C8").Value = "E"

# Row 9
$ws.Range("A9").Value = "Ohh okay :) So the appointment for USER? Right ? :)"
$ws.Range("B9").Value = "NAME_CHECK_N"
$ws.Range("C9").Value = "N"

# Row 10 (new)
$ws.Range("A10").Value = "Ohh sorry! So this appointment for USER? Isn't it ?"
$ws.Range("B10").Value = "NAME_CHECK_N"
$ws.Range("C10").Value = "N"

# Apply the same style (centered alignment) used by columns B/C to the new row 10 cells
$ws.Range("B10:C10").HorizontalAlignment = -4108

# Update selection to match target state
$ws.Range("E11").Select()
